# Habit playbook template - "First changes to update Adoption stuff for Kate"
$p = $ppt.ActivePresentation

# --- Notes Master: refresh the "datetimeFigureOut" date placeholder to the
#     edit date (4/3/2018 -> 9/14/2018). The HeadersFooters/DateAndTime
#     object is the COM surface for this value.
$nm = $p.NotesMaster
$hf = $nm.HeadersFooters
$dt = $hf.DateAndTime
$dt.UseFormat = 0
$dt.Text = "9/14/2018"

# --- Slide 1 edits ---
$s = $p.Slides.Item(1)

# "Content Placeholder 4" (id=5): reword the habit-importance prompt.
$whyShape = $s.Shapes.Item("Content Placeholder 4")
$whyShape.TextFrame.TextRange.Text = "[INSERT WHY THIS HABIT IS IMPORTANT TO DEVELOP]"

# "Content Placeholder 919" (id=920): vertically center the text box content.
$planShape = $s.Shapes.Item("Content Placeholder 919")
$planShape.TextFrame.VerticalAnchor = 3
